$d = $word.ActiveDocument
$d.Content.Find.Execute("77-76=", $true, $false, $false, $false, $false, $true, 1, $false, "61-50=", 2) | Out-Null
$d.Content.Find.Execute("70-39=", $true, $false, $false, $false, $false, $true, 1, $false, "0+86=", 2) | Out-Null
$d.Content.Find.Execute("41+46=", $true, $false, $false, $false, $false, $true, 1, $false, "60-28=", 2) | Out-Null
$d.Content.Find.Execute("62-18=", $true, $false, $false, $false, $false, $true, 1, $false, "58+11=", 2) | Out-Null
$d.Content.Find.Execute("90-15=", $true, $false, $false, $false, $false, $true, 1, $false, "43+6=", 2) | Out-Null
$d.Content.Find.Execute("25-5=", $true, $false, $false, $false, $false, $true, 1, $false, "64-33=", 2) | Out-Null
$d.Content.Find.Execute("52+15=", $true, $false, $false, $false, $false, $true, 1, $false, "99-8=", 2) | Out-Null
$d.Content.Find.Execute("29+18=", $true, $false, $false, $false, $false, $true, 1, $false, "91+2=", 2) | Out-Null
$d.Content.Find.Execute("71+26=", $true, $false, $false, $false, $false, $true, 1, $false, "1+24=", 2) | Out-Null
$d.Content.Find.Execute("10+57=", $true, $false, $false, $false, $false, $true, 1, $false, "71-20=", 2) | Out-Null
$d.Content.Find.Execute("26+20=", $true, $false, $false, $false, $false, $true, 1, $false, "50+37=", 2) | Out-Null
$d.Content.Find.Execute("91-69=", $true, $false, $false, $false, $false, $true, 1, $false, "76-43=", 2) | Out-Null
$d.Content.Find.Execute("90-25=", $true, $false, $false, $false, $false, $true, 1, $false, "37+45=", 2) | Out-Null
$d.Content.Find.Execute("48+3=", $true, $false, $false, $false, $false, $true, 1, $false, "85-28=", 2) | Out-Null
$d.Content.Find.Execute("57-16=", $true, $false, $false, $false, $false, $true, 1, $false, "20+8=", 2) | Out-Null
$d.Content.Find.Execute("95-85=", $true, $false, $false, $false, $false, $true, 1, $false, "71-62=", 2) | Out-Null
$d.Content.Find.Execute("64-54=", $true, $false, $false, $false, $false, $true, 1, $false, "48-19=", 2) | Out-Null
$d.Content.Find.Execute("13+55=", $true, $false, $false, $false, $false, $true, 1, $false, "49-7=", 2) | Out-Null
$d.Content.Find.Execute("71-18=", $true, $false, $false, $false, $false, $true, 1, $false, "14+82=", 2) | Out-Null
$d.Content.Find.Execute("84-13=", $true, $false, $false, $false, $false, $true, 1, $false, "85-25=", 2) | Out-Null
$d.Content.Find.Execute("63+35=", $true, $false, $false, $false, $false, $true, 1, $false, "31-17=", 2) | Out-Null
$d.Content.Find.Execute("55+39=", $true, $false, $false, $false, $false, $true, 1, $false, "38-21=", 2) | Out-Null
$d.Content.Find.Execute("44+36=", $true, $false, $false, $false, $false, $true, 1, $false, "26+63=", 2) | Out-Null
$d.Content.Find.Execute("3+2=", $true, $false, $false, $false, $false, $true, 1, $false, "68-42=", 2) | Out-Null
$d.Content.Find.Execute("20+29=", $true, $false, $false, $false, $false, $true, 1, $false, "14+73=", 2) | Out-Null
$d.Content.Find.Execute("57-0=", $true, $false, $false, $false, $false, $true, 1, $false, "72+12=", 2) | Out-Null
$d.Content.Find.Execute("30-13=", $true, $false, $false, $false, $false, $true, 1, $false, "93-50=", 2) | Out-Null
$d.Content.Find.Execute("19+26=", $true, $false, $false, $false, $false, $true, 1, $false, "27+51=", 2) | Out-Null
$d.Content.Find.Execute("66-6=", $true, $false, $false, $false, $false, $true, 1, $false, "71-8=", 2) | Out-Null
$d.Content.Find.Execute("97-62=", $true, $false, $false, $false, $false, $true, 1, $false, "68-54=", 2) | Out-Null
$d.Content.Find.Execute("72-1=", $true, $false, $false, $false, $false, $true, 1, $false, "54-39=", 2) | Out-Null
$d.Content.Find.Execute("38-25=", $true, $false, $false, $false, $false, $true, 1, $false, "44-28=", 2) | Out-Null
$d.Content.Find.Execute("25+32=", $true, $false, $false, $false, $false, $true, 1, $false, "87-87=", 2) | Out-Null
$d.Content.Find.Execute("61+20=", $true, $false, $false, $false, $false, $true, 1, $false, "34-24=", 2) | Out-Null
$d.Content.Find.Execute("87-64=", $true, $false, $false, $false, $false, $true, 1, $false, "0+63=", 2) | Out-Null
$d.Content.Find.Execute("88-75=", $true, $false, $false, $false, $false, $true, 1, $false, "0+85=", 2) | Out-Null
$d.Content.Find.Execute("93-31=", $true, $false, $false, $false, $false, $true, 1, $false, "52+45=", 2) | Out-Null
$d.Content.Find.Execute("22+19=", $true, $false, $false, $false, $false, $true, 1, $false, "60+29=", 2) | Out-Null
$d.Content.Find.Execute("8+60=", $true, $false, $false, $false, $false, $true, 1, $false, "36+42=", 2) | Out-Null
$d.Content.Find.Execute("58+30=", $true, $false, $false, $false, $false, $true, 1, $false, "58-45=", 2) | Out-Null
$d.Content.Find.Execute("35+57=", $true, $false, $false, $false, $false, $true, 1, $false, "54+20=", 2) | Out-Null
$d.Content.Find.Execute("67-14=", $true, $false, $false, $false, $false, $true, 1, $false, "8+19=", 2) | Out-Null
$d.Content.Find.Execute("35-4=", $true, $false, $false, $false, $false, $true, 1, $false, "18-10=", 2) | Out-Null
$d.Content.Find.Execute("26-21=", $true, $false, $false, $false, $false, $true, 1, $false, "31+28=", 2) | Out-Null
$d.Content.Find.Execute("84-62=", $true, $false, $false, $false, $false, $true, 1, $false, "95-44=", 2) | Out-Null
$d.Content.Find.Execute("63-33=", $true, $false, $false, $false, $false, $true, 1, $false, "99-97=", 2) | Out-Null
$d.Content.Find.Execute("30+9=", $true, $false, $false, $false, $false, $true, 1, $false, "46+45=", 2) | Out-Null
$d.Content.Find.Execute("18+6=", $true, $false, $false, $false, $false, $true, 1, $false, "36+44=", 2) | Out-Null
$d.Content.Find.Execute("73-36=", $true, $false, $false, $false, $false, $true, 1, $false, "24+72=", 2) | Out-Null
$d.Content.Find.Execute("26+0=", $true, $false, $false, $false, $false, $true, 1, $false, "86-64=", 2) | Out-Null
$d.Content.Find.Execute("0+47=", $true, $false, $false, $false, $false, $true, 1, $false, "58-29=", 2) | Out-Null
$d.Content.Find.Execute("4+10=", $true, $false, $false, $false, $false, $true, 1, $false, "77-52=", 2) | Out-Null
$d.Content.Find.Execute("86-75=", $true, $false, $false, $false, $false, $true, 1, $false, "59+14=", 2) | Out-Null
$d.Content.Find.Execute("47-24=", $true, $false, $false, $false, $false, $true, 1, $false, "36+33=", 2) | Out-Null
$d.Content.Find.Execute("28+21=", $true, $false, $false, $false, $false, $true, 1, $false, "15-9=", 2) | Out-Null
$d.Content.Find.Execute("81-4=", $true, $false, $false, $false, $false, $true, 1, $false, "91-0=", 2) | Out-Null
$d.Content.Find.Execute("43-22=", $true, $false, $false, $false, $false, $true, 1, $false, "15-6=", 2) | Out-Null
$d.Content.Find.Execute("24+57=", $true, $false, $false, $false, $false, $true, 1, $false, "96-37=", 2) | Out-Null
$d.Content.Find.Execute("51-2=", $true, $false, $false, $false, $false, $true, 1, $false, "62+1=", 2) | Out-Null
$d.Content.Find.Execute("80-73=", $true, $false, $false, $false, $false, $true, 1, $false, "36+10=", 2) | Out-Null
$d.Content.Find.Execute("64-46=", $true, $false, $false, $false, $false, $true, 1, $false, "79-17=", 2) | Out-Null
$d.Content.Find.Execute("12+24=", $true, $false, $false, $false, $false, $true, 1, $false, "43+44=", 2) | Out-Null
$d.Content.Find.Execute("73-9=", $true, $false, $false, $false, $false, $true, 1, $false, "55+9=", 2) | Out-Null
$d.Content.Find.Execute("87-41=", $true, $false, $false, $false, $false, $true, 1, $false, "59-24=", 2) | Out-Null
$d.Content.Find.Execute("79-27=", $true, $false, $false, $false, $false, $true, 1, $false, "12+16=", 2) | Out-Null
$d.Content.Find.Execute("81-80=", $true, $false, $false, $false, $false, $true, 1, $false, "82-27=", 2) | Out-Null
$d.Content.Find.Execute("22+12=", $true, $false, $false, $false, $false, $true, 1, $false, "88-7=", 2) | Out-Null
$d.Content.Find.Execute("97-52=", $true, $false, $false, $false, $false, $true, 1, $false, "19-2=", 2) | Out-Null
$d.Content.Find.Execute("78-6=", $true, $false, $false, $false, $false, $true, 1, $false, "18-9=", 2) | Out-Null
$d.Content.Find.Execute("32+37=", $true, $false, $false, $false, $false, $true, 1, $false, "6+16=", 2) | Out-Null
$d.Content.Find.Execute("20+66=", $true, $false, $false, $false, $false, $true, 1, $false, "55+0=", 2) | Out-Null
$d.Content.Find.Execute("26+32=", $true, $false, $false, $false, $false, $true, 1, $false, "52-40=", 2) | Out-Null
$d.Content.Find.Execute("13+81=", $true, $false, $false, $false, $false, $true, 1, $false, "82-24=", 2) | Out-Null
$d.Content.Find.Execute("16+23=", $true, $false, $false, $false, $false, $true, 1, $false, "55+22=", 2) | Out-Null
$d.Content.Find.Execute("38+14=", $true, $false, $false, $false, $false, $true, 1, $false, "19+23=", 2) | Out-Null
$d.Content.Find.Execute("82-72=", $true, $false, $false, $false, $false, $true, 1, $false, "10+76=", 2) | Out-Null
$d.Content.Find.Execute("94-71=", $true, $false, $false, $false, $false, $true, 1, $false, "79-33=", 2) | Out-Null
$d.Content.Find.Execute("23+65=", $true, $false, $false, $false, $false, $true, 1, $false, "86-1=", 2) | Out-Null
$d.Content.Find.Execute("81-36=", $true, $false, $false, $false, $false, $true, 1, $false, "74-60=", 2) | Out-Null
$d.Content.Find.Execute("13+37=", $true, $false, $false, $false, $false, $true, 1, $false, "5+86=", 2) | Out-Null
$d.Content.Find.Execute("84+11=", $true, $false, $false, $false, $false, $true, 1, $false, "33+60=", 2) | Out-Null
$d.Content.Find.Execute("5+42=", $true, $false, $false, $false, $false, $true, 1, $false, "89-29=", 2) | Out-Null
$d.Content.Find.Execute("58+4=", $true, $false, $false, $false, $false, $true, 1, $false, "19+69=", 2) | Out-Null
$d.Content.Find.Execute("53-37=", $true, $false, $false, $false, $false, $true, 1, $false, "86+0=", 2) | Out-Null
$d.Content.Find.Execute("24+39=", $true, $false, $false, $false, $false, $true, 1, $false, "89-10=", 2) | Out-Null
$d.Content.Find.Execute("13+41=", $true, $false, $false, $false, $false, $true, 1, $false, "67-65=", 2) | Out-Null
$d.Content.Find.Execute("34+7=", $true, $false, $false, $false, $false, $true, 1, $false, "49+0=", 2) | Out-Null
$d.Content.Find.Execute("15+14=", $true, $false, $false, $false, $false, $true, 1, $false, "34-0=", 2) | Out-Null
$d.Content.Find.Execute("0+6=", $true, $false, $false, $false, $false, $true, 1, $false, "67-63=", 2) | Out-Null
$d.Content.Find.Execute("34+36=", $true, $false, $false, $false, $false, $true, 1, $false, "65-42=", 2) | Out-Null
$d.Content.Find.Execute("81-40=", $true, $false, $false, $false, $false, $true, 1, $false, "94-60=", 2) | Out-Null
$d.Content.Find.Execute("92-3=", $true, $false, $false, $false, $false, $true, 1, $false, "51+0=", 2) | Out-Null
$d.Content.Find.Execute("76-27=", $true, $false, $false, $false, $false, $true, 1, $false, "58+25=", 2) | Out-Null
$d.Content.Find.Execute("43-20=", $true, $false, $false, $false, $false, $true, 1, $false, "35+41=", 2) | Out-Null
$d.Content.Find.Execute("69-52=", $true, $false, $false, $false, $false, $true, 1, $false, "73-69=", 2) | Out-Null
$d.Content.Find.Execute("6+62=", $true, $false, $false, $false, $false, $true, 1, $false, "4+61=", 2) | Out-Null
$d.Content.Find.Execute("17+57=", $true, $false, $false, $false, $false, $true, 1, $false, "21+20=", 2) | Out-Null
$d.Content.Find.Execute("53-26=", $true, $false, $false, $false, $false, $true, 1, $false, "59+30=", 2) | Out-Null
$d.Content.Find.Execute("96-32=", $true, $false, $false, $false, $false, $true, 1, $false, "75-68=", 2) | Out-Null
$d.Content.Find.Execute("15+33=", $true, $false, $false, $false, $false, $true, 1, $false, "84+14=", 2) | Out-Null
